# Lattice multiplication exercises: replace the multiplication problems
# in all 15 table cells with the new set of values, preserving the
# existing run formatting (sz=32) and line-break structure.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellData = @(
    @{ Row = 1; Col = 1; Title = '12 x 60'; Row2 = '  6    0'; L1 = '1|    |'; L2 = '2|    |' }
    @{ Row = 1; Col = 2; Title = '47 x 36'; Row2 = '  3    6'; L1 = '4|    |'; L2 = '7|    |' }
    @{ Row = 1; Col = 3; Title = '63 x 93'; Row2 = '  9    3'; L1 = '6|    |'; L2 = '3|    |' }
    @{ Row = 2; Col = 1; Title = '19 x 12'; Row2 = '  1    2'; L1 = '1|    |'; L2 = '9|    |' }
    @{ Row = 2; Col = 2; Title = '78 x 55'; Row2 = '  5    5'; L1 = '7|    |'; L2 = '8|    |' }
    @{ Row = 2; Col = 3; Title = '93 x 69'; Row2 = '  6    9'; L1 = '9|    |'; L2 = '3|    |' }
    @{ Row = 3; Col = 1; Title = '41 x 82'; Row2 = '  8    2'; L1 = '4|    |'; L2 = '1|    |' }
    @{ Row = 3; Col = 2; Title = '95 x 71'; Row2 = '  7    1'; L1 = '9|    |'; L2 = '5|    |' }
    @{ Row = 3; Col = 3; Title = '23 x 42'; Row2 = '  4    2'; L1 = '2|    |'; L2 = '3|    |' }
    @{ Row = 4; Col = 1; Title = '15 x 76'; Row2 = '  7    6'; L1 = '1|    |'; L2 = '5|    |' }
    @{ Row = 4; Col = 2; Title = '63 x 12'; Row2 = '  1    2'; L1 = '6|    |'; L2 = '3|    |' }
    @{ Row = 4; Col = 3; Title = '37 x 10'; Row2 = '  1    0'; L1 = '3|    |'; L2 = '7|    |' }
    @{ Row = 5; Col = 1; Title = '50 x 42'; Row2 = '  4    2'; L1 = '5|    |'; L2 = '0|    |' }
    @{ Row = 5; Col = 2; Title = '98 x 83'; Row2 = '  8    3'; L1 = '9|    |'; L2 = '8|    |' }
    @{ Row = 5; Col = 3; Title = '88 x 54'; Row2 = '  5    4'; L1 = '8|    |'; L2 = '8|    |' }
)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

foreach ($item in $cellData) {
    $cell = $t.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    $xml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" + 
            "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" + 
            "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" + 
            "<pkg:xmlData>" + 
            "<w:document xmlns:w=`"$wNs`">" + 
            "<w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr>" + 
            "<w:t>$($item.Title)</w:t><w:br/>" + 
            "<w:t xml:space=`"preserve`">$($item.Row2)</w:t><w:br/>" + 
            "<w:t xml:space=`"preserve`">  ----</w:t><w:br/>" + 
            "<w:t>$($item.L1)</w:t><w:br/>" + 
            "<w:t>$($item.L2)</w:t>" + 
            "</w:r></w:p></w:body></w:document>" + 
            "</pkg:xmlData></pkg:part></pkg:package>"
    $rng.InsertXML($xml) | Out-Null
}

Write-Output "Updated $($cellData.Count) cells"